# Bugfixed evaluation and simulated rt_data for components.
#
# The forecast vector table (date_of_forecast, y_0, y_0_forecast, y_1,
# y_1_forecast) is re-derived: every existing row shifts down by one
# forecast period and gains a y_0_forecast value (column C) that used to
# be missing for several rows, and one new trailing row (53) is appended
# for the latest forecast date/period.
#
# Simplest faithful reproduction: rewrite the whole A2:E53 block with the
# corrected values, then restore the date-column cell style (s="2", the
# bold/centered/bordered custom date format already used by A2:A52) onto
# the newly created A53 cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(39400,2007,1.144978573787081,2008,1.854677717601594),
    @(39583,2008,1.981300392979479,2009,2.228397109637226),
    @(39765,2008,1.381024225294869,2009,1.929644353290927),
    @(39948,2009,-0.4626567965509865,2010,-0.2018858887078534),
    @(40130,2009,-0.8792832172735965,2010,-0.10374497415091),
    @(40310,2010,-0.4861108058567654,2011,0.2108047537406676),
    @(40494,2010,0.9337833426867448,2011,1.899713704008654),
    @(40676,2011,2.490556233265284,2012,0.9684279156219944),
    @(40862,2011,2.791140000794279,2012,1.722110645261998),
    @(41044,2012,1.175152734113438,2013,1.69818237209749),
    @(41228,2012,0.4451370000809973,2013,0.6480763427741953),
    @(41409,2013,-0.2674335569108899,2014,0.679606031449409),
    @(41592,2013,0.2545814083968478,2014,1.223618887196531),
    @(41774,2014,2.038609866767938,2015,1.325176859452326),
    @(41957,2014,1.297015177357297,2015,0.309136676902555),
    @(42137,2015,1.566972126803345,2016,1.522808462763714),
    @(42321,2015,1.365576377841027,2016,1.582979977679533),
    @(42503,2016,2.246337373619012,2017,1.693557061600992),
    @(42689,2016,2.204449574611278,2017,1.401113624217043),
    @(42867,2017,1.777150434343522,2018,1.905564797014625),
    @(43053,2017,2.18621550610123,2018,2.357704431248409),
    @(43145,2018,1.953801996161997,2019,2.285039356133778),
    @(43235,2018,1.741137453897301,2019,2.082477074609068),
    @(43326,2018,1.963716944695593,2019,2.26303158209491),
    @(43418,2018,1.911050033324102,2019,2.362613045431528),
    @(43510,2019,1.83152322656317,2020,2.06914762522894),
    @(43600,2019,1.562095320687407,2020,1.845103901518885),
    @(43691,2019,1.657291696732699,2020,1.670042545405925),
    @(43783,2019,1.457852003181337,2020,0.7756897792100315),
    @(43875,2020,0.9716673922241847,2021,1.221567067111873),
    @(43966,2020,-0.8985012482809474,2021,0.02570757229449772),
    @(44068,2020,-3.357986809108382,2021,-1.069210638149831),
    @(44159,2020,-3.258619210312885,2021,-2.035749133083187),
    @(44251,2021,-2.792286268255573,2022,-2.462477360333626),
    @(44341,2021,-0.7241284594088016,2022,-0.1006764062950749),
    @(44432,2021,-0.1388955462784724,2022,1.381312451245953),
    @(44525,2021,0.4255262881966981,2022,2.107044281672787),
    @(44617,2022,2.610016230858681,2023,0.4182551040150306),
    @(44706,2022,3.72986100763808,2023,1.357201584213352),
    @(44798,2022,3.463320568938566,2023,1.409962909690821),
    @(44890,2022,3.293290997728171,2023,1.50842747477733),
    @(44981,2023,0.2608494586689769,2024,1.34681835462207),
    @(45071,2023,-0.3358560228777674,2024,1.155533515694374),
    @(45163,2023,-0.2871245688614854,2024,1.278647349573703),
    @(45254,2023,-0.2814561130375925,2024,0.06837453075889677),
    @(45345,2024,-0.4741729829813357,2025,0.061310808297943),
    @(45436,2024,-1.054811008161194,2025,-0.009717596728553435),
    @(45534,2024,-0.9636841177852018,2025,-0.2108962101654677),
    @(45618,2024,-0.6470065423293869,2025,1.000784331770643),
    @(45713,2025,2.411585640991931,2026,0.3354683421588955),
    @(45800,2025,1.593155398714674,2026,0.6473947787101642),
    @(45891,2025,2.053865394798304,2026,0.8520203638610502)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $arr = New-Object 'object[,]' 1,5
    for ($j = 0; $j -lt 5; $j++) {
        $arr[0,$j] = $row[$j]
    }
    $ws.Range("A" + $r + ":E" + $r).Value2 = $arr
}

# The table grew by one row (52 -> 53 data rows, A1:E52 -> A1:E53). The
# brand-new row 53 has no style yet, so clone the date-column style
# (s="2") from the row directly above it onto A53. -4122 = xlPasteFormats.
$xlPasteFormats = -4122
$ws.Range("A52").Copy()
$ws.Range("A53").PasteSpecial($xlPasteFormats)
